# FAST_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer from 2021-03-17 to 2021-03-18
#  - refresh Weight (col D) and Percent Change (col E) figures for rows 2-10
#
# The sheet ships protected (no UI password known to us), so we have to lift
# protection before writing and then restore it so the sheet goes back to a
# protected state once we're done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# --- Disclaimer text: only the date changes -------------------------------
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."
# Re-fit the row so the wrapped-text edit doesn't leave a stray explicit
# row height behind (the row was auto-height before the edit).
$ws.Rows.Item(13).AutoFit()

# --- Weight (D) / Percent Change (E) refresh -------------------------------
$ws.Range("D2").Value = 0.1092499541036314
$ws.Range("E2").Value = -0.04656414668243902

$ws.Range("D3").Value = 0.1049329337173569
$ws.Range("E3").Value = -0.03370233702337011

$ws.Range("D4").Value = 0.1149849110796559
$ws.Range("E4").Value = -0.01460093511606919

$ws.Range("D5").Value = 0.1365277578718573
$ws.Range("E5").Value = -0.01363829667993022

$ws.Range("D6").Value = 0.1301989785170908
$ws.Range("E6").Value = -0.007161648640748197

$ws.Range("D7").Value = 0.1452381411843814
$ws.Range("E7").Value = -0.01100583787922293

$ws.Range("D8").Value = 0.1282199235966322
$ws.Range("E8").Value = -0.02151973879489455

$ws.Range("D9").Value = 0.1306473999293941
$ws.Range("E9").Value = -0.02353564279471398

$ws.Range("E10").Value = -0.02052954585200351

$ws.Protect()
